$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6: change the table's style (Table Design gallery) from the
#    deck's local custom style to the built-in "Medium Style 2 - Accent 1".
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{16CB49F4-266B-47C4-9DEC-DAF4CC452A3D}")

# ---------------------------------------------------------------------------
# 2) Design tab: switch the presentation's theme from "Integral" to the
#    plain default "Office Theme" (same font/format scheme, new palette).
# ---------------------------------------------------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0x000000   # Text/Background - Dark 1   -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # Text/Background - Light 1  -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # Text/Background - Dark 2   -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # Text/Background - Light 2  -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # Accent 1                    -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # Accent 2                    -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # Accent 3                    -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # Accent 4                    -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # Accent 5                    -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # Accent 6                    -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # Hyperlink                   -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # Followed Hyperlink          -> 954F72
